$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column D ("Tipo"), shifting it to E
$ws.Range("D1").EntireColumn.Insert()

# New header for the inserted column, matching the other header cells' formatting
$ws.Range("D1").Value = "MAE"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1

# Updated MSE (B) and R2 (C) values, plus new MAE (D) values
$ws.Range("B2").Value = 0.5022389224929714
$ws.Range("C2").Value = 0.9899993564879116
$ws.Range("D2").Value = 0.5794856473835275

$ws.Range("B3").Value = 0.2364677015639231
$ws.Range("C3").Value = 0.9953782328262166
$ws.Range("D3").Value = 0.3872941361072821

$ws.Range("B4").Value = 0.2812201673132656
$ws.Range("C4").Value = 0.9945899493137279
$ws.Range("D4").Value = 0.4300862106884557

$ws.Range("B5").Value = 0.4130513438249112
$ws.Range("C5").Value = 0.9918551930174837
$ws.Range("D5").Value = 0.4941922081108229
